# Trade #194 closed at 2026-02-17 22:08:02 - unknown UNKNOWN +0.000%
#
# This script updates the live trading workbook:
#  1) Summary sheet aggregate metrics
#  2) Strategy Status sheet for MarketMaking strategy row
#  3) All Trades sheet: closes trade #222 (row 223) and appends 2 new open trades
#  4) volatility_scorer sheet: appends the newly opened volatility_scorer trade
#  5) MarketMaking sheet: closes trade #222 (row 190) and appends the new MarketMaking trade

$wb = $excel.ActiveWorkbook

# Helper: set a text value into a cell while preventing Excel's automatic
# date/time literal conversion (e.g. "2026-02-17" becoming a date serial).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1399.83   # Current Capital
$wsSummary.Range("B4").Value = -0.39     # Total P&L $
$wsSummary.Range("B5").Value = -0.04     # Total P&L %
$wsSummary.Range("B6").Value = 222       # Total Trades
$wsSummary.Range("B8").Value = 98        # Losing Trades
$wsSummary.Range("B9").Value = 38.29     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 99.83
$wsStatus.Range("D5").Value = 189
$wsStatus.Range("E5").Value = -0.5
$wsStatus.Range("F5").Value = -0.17
$wsStatus.Range("G5").Value = 37.57

# ---------------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close trade #222 (row 223)
$wsAll.Cells.Item(223, 7).Value = 0.33              # G223 Exit Price
$wsAll.Cells.Item(223, 8).Value = "CLOSED"          # H223 Status
$wsAll.Cells.Item(223, 9).Value = -10.8108          # I223 P&L %
$wsAll.Cells.Item(223, 10).Value = -0.04            # J223 P&L $
$wsAll.Cells.Item(223, 11).Value = 99.83            # K223 Capital After
$wsAll.Cells.Item(223, 12).Value = "early_exit"     # L223 Exit Reason
$wsAll.Cells.Item(223, 13).Value = 0.18             # M223 Duration (min)

# New row 256: trade #255 (volatility_scorer, OPEN)
Set-TextValue $wsAll.Cells.Item(256, 2) "2026-02-17"
Set-TextValue $wsAll.Cells.Item(256, 3) "22:07:54"
$wsAll.Cells.Item(256, 1).Value = 255
$wsAll.Cells.Item(256, 4).Value = "volatility_scorer"
$wsAll.Cells.Item(256, 5).Value = "NEUTRAL"
$wsAll.Cells.Item(256, 6).Value = 0.37
$wsAll.Cells.Item(256, 8).Value = "OPEN"
$wsAll.Cells.Item(256, 9).Value = 0
$wsAll.Cells.Item(256, 10).Value = 0
$wsAll.Cells.Item(256, 11).Value = 100
$wsAll.Cells.Item(256, 13).Value = 0
$wsAll.Cells.Item(256, 14).Value = 0
$wsAll.Cells.Item(256, 15).Value = 0
$wsAll.Cells.Item(256, 16).Value = 0.85
$wsAll.Cells.Item(256, 17).Value = "Low vol market (score: inf) - ideal for market making"

# New row 257: trade #256 (MarketMaking, OPEN)
Set-TextValue $wsAll.Cells.Item(257, 2) "2026-02-17"
Set-TextValue $wsAll.Cells.Item(257, 3) "22:07:55"
$wsAll.Cells.Item(257, 1).Value = 256
$wsAll.Cells.Item(257, 4).Value = "MarketMaking"
$wsAll.Cells.Item(257, 5).Value = "UP"
$wsAll.Cells.Item(257, 6).Value = 0.64
$wsAll.Cells.Item(257, 8).Value = "OPEN"
$wsAll.Cells.Item(257, 9).Value = 0
$wsAll.Cells.Item(257, 10).Value = 0
$wsAll.Cells.Item(257, 11).Value = 99.86837680355363
$wsAll.Cells.Item(257, 13).Value = 0
$wsAll.Cells.Item(257, 14).Value = 0
$wsAll.Cells.Item(257, 15).Value = 0
$wsAll.Cells.Item(257, 16).Value = 0.6
$wsAll.Cells.Item(257, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# 4) volatility_scorer sheet - append trade #255 (row 10)
#    Column layout: A Trade#, B Date, C Time, D Strategy, E Side, F Entry,
#    G Exit, H Status, I P&L%, J P&L$, K Capital After, L Entry Slip,
#    M Exit Slip, N Confidence, O Entry Reason, P Exit Reason, Q Duration
# ---------------------------------------------------------------------------
$wsVol = $wb.Worksheets.Item("volatility_scorer")
Set-TextValue $wsVol.Cells.Item(10, 2) "2026-02-17"
Set-TextValue $wsVol.Cells.Item(10, 3) "22:07:54"
$wsVol.Cells.Item(10, 1).Value = 255
$wsVol.Cells.Item(10, 4).Value = "volatility_scorer"
$wsVol.Cells.Item(10, 5).Value = "NEUTRAL"
$wsVol.Cells.Item(10, 6).Value = 0.37
$wsVol.Cells.Item(10, 8).Value = "OPEN"
$wsVol.Cells.Item(10, 9).Value = 0
$wsVol.Cells.Item(10, 10).Value = 0
$wsVol.Cells.Item(10, 11).Value = 100
$wsVol.Cells.Item(10, 12).Value = 0
$wsVol.Cells.Item(10, 13).Value = 0
$wsVol.Cells.Item(10, 14).Value = 0.85
$wsVol.Cells.Item(10, 15).Value = "Low vol market (score: inf) - ideal for market making"
$wsVol.Cells.Item(10, 17).Value = 0

# ---------------------------------------------------------------------------
# 5) MarketMaking sheet
#    Column layout: A Trade#, B Date, C Time, D Strategy, E Side, F Entry,
#    G Exit, H Status, I P&L%, J P&L$, K Capital After, L Entry Slip,
#    M Exit Slip, N Confidence, O Entry Reason, P Exit Reason, Q Duration
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Close trade #222 (row 190)
$wsMM.Cells.Item(190, 7).Value = 0.33                 # G190 Exit Price
$wsMM.Cells.Item(190, 8).Value = "CLOSED"             # H190 Status
$wsMM.Cells.Item(190, 9).Value = -10.8108             # I190 P&L %
$wsMM.Cells.Item(190, 10).Value = -0.04               # J190 P&L $
$wsMM.Cells.Item(190, 11).Value = 99.83               # K190 Capital After
$wsMM.Cells.Item(190, 16).Value = "early_exit"        # P190 Exit Reason
$wsMM.Cells.Item(190, 17).Value = 0.18                # Q190 Duration (min)

# New row 215: trade #256 (MarketMaking, OPEN)
Set-TextValue $wsMM.Cells.Item(215, 2) "2026-02-17"
Set-TextValue $wsMM.Cells.Item(215, 3) "22:07:55"
$wsMM.Cells.Item(215, 1).Value = 256
$wsMM.Cells.Item(215, 4).Value = "MarketMaking"
$wsMM.Cells.Item(215, 5).Value = "UP"
$wsMM.Cells.Item(215, 6).Value = 0.64
$wsMM.Cells.Item(215, 8).Value = "OPEN"
$wsMM.Cells.Item(215, 9).Value = 0
$wsMM.Cells.Item(215, 10).Value = 0
$wsMM.Cells.Item(215, 11).Value = 99.86837680355363
$wsMM.Cells.Item(215, 12).Value = 0
$wsMM.Cells.Item(215, 13).Value = 0
$wsMM.Cells.Item(215, 14).Value = 0.6
$wsMM.Cells.Item(215, 15).Value = "Normal spread capture: 19600 bps"
$wsMM.Cells.Item(215, 17).Value = 0
